$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.075.55"
$ws.Range("E2").Value = "  +2.76%  "
$ws.Range("D3").Value = "3.204.31"
$ws.Range("E3").Value = "  +1.53%  "
$ws.Range("D5").Value = "'537.69"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").Value = "'145.86"
$ws.Range("E6").Value = "  +4.26%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'0.532"
$ws.Range("E8").Value = "  +3.04%  "
$ws.Range("E9").Value = "  +0.61%  "
$ws.Range("D10").Value = "'0.113"
$ws.Range("E10").Value = "  +3.57%  "
$ws.Range("D11").Value = "'0.435"
$ws.Range("E11").Value = "  +3.00%  "
$ws.Range("D12").Value = "3.760.00"
$ws.Range("E12").Value = "  +1.40%  "
$ws.Range("E13").Value = "  -0.99%  "
$ws.Range("D14").Value = "'26.18"
$ws.Range("E14").Value = "  +0.21%  "
$ws.Range("E15").Value = "  +2.96%  "
$ws.Range("D16").Value = "60.179.42"
$ws.Range("E16").Value = "  +2.83%  "
$ws.Range("D17").Value = "3.214.04"
$ws.Range("E17").Value = "  +1.09%  "
$ws.Range("D18").Value = "'6.26"
$ws.Range("E18").Value = "  +0.67%  "
$ws.Range("D19").Value = "'13.24"
$ws.Range("E19").Value = "  +1.37%  "
$ws.Range("D20").Value = "'8.33"
$ws.Range("E20").Value = "  +0.98%  "
$ws.Range("D21").Value = "'379.75"
$ws.Range("E21").Value = "  +0.56%  "
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("E23").Value = "  +1.87%  "
$ws.Range("D24").Value = "'70.16"
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("D25").Value = "'8.95"
$ws.Range("E25").Value = "  +10.08%  "
$ws.Range("E26").Value = "  +1.49%  "
$ws.Range("E27").Value = "  +1.29%  "
$ws.Range("D28").Value = "0.0₃0904"
$ws.Range("E28").Value = "  +3.49%  "
$ws.Range("D29").Value = "'6.22"
$ws.Range("E29").Value = "  +1.22%  "
$ws.Range("E30").Value = "  +0.76%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "'5.44"
$ws.Range("E31").Value = "  +5.03%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "'22.34"
$ws.Range("E32").Value = "  +2.09%  "
$ws.Range("E33").Value = "  +3.74%  "
$ws.Range("D34").Value = "'6.67"
$ws.Range("E34").Value = "  +6.65%  "
$ws.Range("D35").Value = "'156.94"
$ws.Range("E35").Value = "  -2.45%  "
$ws.Range("E36").Value = "  -0.86%  "
$ws.Range("D37").Value = "2.797.93"
$ws.Range("E37").Value = "  +6.14%  "
$ws.Range("D38").Value = "'25.64"
$ws.Range("E38").Value = "  +0.79%  "
$ws.Range("D39").Value = "'0.0705"
$ws.Range("E39").Value = "  +3.30%  "
$ws.Range("E40").Value = "  +0.45%  "
$ws.Range("E41").Value = "  +1.02%  "
$ws.Range("D42").Value = "'39.84"
$ws.Range("E42").Value = "  +3.05%  "
$ws.Range("D43").Value = "'0.0295"
$ws.Range("E43").Value = "  +5.01%  "
$ws.Range("D44").Value = "'0.718"
$ws.Range("E44").Value = "  +1.63%  "
$ws.Range("E45").Value = "  +2.43%  "
$ws.Range("D46").Value = "3.249.32"
$ws.Range("E46").Value = "  +1.51%  "
$ws.Range("E47").Value = "  +2.20%  "
$ws.Range("E48").Value = "  -0.50%  "
$ws.Range("D49").Value = "'0.808"
$ws.Range("E49").Value = "  +6.57%  "
$ws.Range("D50").Value = "'20.71"
$ws.Range("E50").Value = "  +2.14%  "
$ws.Range("E51").Value = "  +0.01%  "
